$d = $word.ActiveDocument

function Wrap-PkgXml($bodyInner) {
  return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Locate the existing "Please note: ..." paragraph via Find.
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute("Please note:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $target.Find.Found) {
    throw "Could not find 'Please note:' paragraph"
}
$targetPara = $target.Paragraphs.Item(1)
$r = $targetPara.Range
$r.Collapse(1)  # wdCollapseStart

# InsertXML on a collapsed range replaces the paragraph at the insertion
# point with the *last* <w:p> supplied, and inserts any preceding <w:p>
# elements as brand-new paragraphs before it. So we must re-supply the
# (modified) "Please note" paragraph as the final paragraph in the payload.

$newPara = '<w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr><w:bookmarkStart w:id="5" w:name="_Hlk150165830"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>The Claimant (or in the event they act in person and the Defendant is represented, the Defendant) must bring to court for the start of the trial a paper copy of the electronic trial bundle for use by witnesses. A failure to do so may result in the imposition of sanctions.</w:t></w:r></w:p>'

$bookmarkEndAndBlank = '<w:bookmarkEnd w:id="5"/><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'

$pleaseNotePara = '<w:p><w:pPr><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Please note: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Cases are listed in accordance with local hearing arrangements determined by the Judiciary and implemented by the court staff. Every effort is made to ensure that hearings start at the time specified. However, listing practices or other factors may mean that you experience a delay, an adjournment at short notice or your case may be released to a different court hearing centre, in which case you will be notified.</w:t></w:r></w:p>'

$payload = $newPara + $bookmarkEndAndBlank + $pleaseNotePara
$xml = Wrap-PkgXml $payload
$r.InsertXML($xml)
